# Update countries & provincias Spain
# Applies the 12-Jul-2020 18:08 -> 19:25 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 19:25"

# --- Updated per-country figures -------------------------------------------------
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes
# Rows 74/75 and 124/125 also swap position in the ranking (country names
# exchange rows) because the updated "Casos totales" value changes their
# sort order: Kenia overtakes Australia, Sierra Leona overtakes Cabo Verde.

$rows = @(
    @{ Row = 4;   Name = $null;           Vals = @(3387307, 31661, 1505671, 1744015, 0, 219, 137621) },  # Estados Unidos
    @{ Row = 6;   Name = $null;           Vals = @(878513,  28155, 553843,  301491,  0, 492, 23179)  },  # India
    @{ Row = 18;  Name = $null;           Vals = @(212993,  1012,  194515,  13115,   0, 19,  5363)   },  # Turquia
    @{ Row = 19;  Name = $null;           Vals = @(199914,  102,   184500,  6280,    0, 0,   9134)   },  # Alemania
    @{ Row = 37;  Name = $null;           Vals = @(54854,   401,   45140,   9381,    0, 2,   333)    },  # Emiratos Arabes Unidos
    @{ Row = 46;  Name = $null;           Vals = @(38670,   1206,  19008,   19300,   0, 8,   362)    },  # Israel
    @{ Row = 56;  Name = $null;           Vals = @(25628,   17,    23364,   518,     0, 0,   1746)   },  # Irlanda
    @{ Row = 65;  Name = $null;           Vals = @(15745,   203,   12283,   3212,    0, 5,   250)    },  # Marruecos
    @{ Row = 74;  Name = "Kenia";         Vals = @(10105,   379,   2881,    7039,    0, 1,   185)    },  # was Australia
    @{ Row = 75;  Name = "Australia";     Vals = @(9796,    243,   7727,    1961,    0, 1,   108)    },  # was Kenia
    @{ Row = 83;  Name = $null;           Vals = @(7560,    158,   2430,    5003,    0, 3,   127)    },  # Etiopia
    @{ Row = 89;  Name = $null;           Vals = @(6552,    46,    5228,    1269,    0, 0,   55)     },  # Tayikistan
    @{ Row = 105; Name = $null;           Vals = @(3059,    8,     1264,    1702,    0, 1,   93)     },  # Somalia
    @{ Row = 110; Name = $null;           Vals = @(2615,    104,   1981,    623,     0, 0,   11)     },  # Sri Lanka
    @{ Row = 111; Name = $null;           Vals = @(2426,    6,     2258,    81,      0, 0,   87)     },  # Cuba
    @{ Row = 124; Name = "Sierra Leona";  Vals = @(1635,    17,    1154,    418,     0, 0,   63)     },  # was Cabo Verde
    @{ Row = 125; Name = "Cabo Verde";    Vals = @(1623,    0,     748,     856,     0, 0,   19)     },  # was Sierra Leona
    @{ Row = 137; Name = $null;           Vals = @(1157,    22,    364,     784,     0, 0,   9)      },  # Mozambique
    @{ Row = 156; Name = $null;           Vals = @(593,     16,    518,     118,     0, 0,   3)      }   # Reunion
)

$cols = @("B", "C", "D", "E", "F", "G", "H")

foreach ($r in $rows) {
    if ($r.Name) {
        $ws.Range("A" + $r.Row).Value = $r.Name
    }
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r.Row).Value = $r.Vals[$i]
    }
}
